$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$newRows = @(
    @(234, 44308, 4, 37, 153.3869496724981),
    @(235, 44309, 5, 32, 132.6589835005389),
    @(236, 44310, 5, 25, 103.639830859796),
    @(237, 44311, 4, 25, 103.639830859796),
    @(238, 44312, 4, 24, 99.4942376254042)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $dateSerial = $row[1]
    $col2 = $row[2]
    $col3 = $row[3]
    $col4 = $row[4]

    # Copy the formatting of the last existing row (233) into the new row so that
    # the date column keeps its style (border/alignment/number format).
    $ws.Range("A233:D233").Copy()
    $ws.Range("A$r:D$r").PasteSpecial(-4122)  # xlPasteFormats

    $ws.Cells.Item($r, 1).Value = $dateSerial
    $ws.Cells.Item($r, 2).Value = $col2
    $ws.Cells.Item($r, 3).Value = $col3
    $ws.Cells.Item($r, 4).Value = $col4
}

$excel.CutCopyMode = 0
